# Updated cryptos list on Sun Mar 17 21:20:32 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    # Writing via an apostrophe prefix forces Excel to keep the value as
    # literal text even when it looks like a number (e.g. "197.73"),
    # matching how the source data (inline strings) is represented.
    # Resetting the style back to Normal afterwards clears the
    # quote-prefix cell style so no stray style index is left on the cell.
    $ws.Range($rangeAddress).Value = "'" + $value
    $ws.Range($rangeAddress).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "68.327.16"
Set-TextValue "E2" "  +2.12%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.636.97"
Set-TextValue "E3" "  +1.03%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.04%  "

# Row 5 - Solana
Set-TextValue "D5" "197.73"
Set-TextValue "E5" "  +8.83%  "

# Row 6 - BNB
Set-TextValue "D6" "577.80"
Set-TextValue "E6" "  -1.46%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.631.89"
Set-TextValue "E7" "  +1.12%  "

# Row 8 - XRP
Set-TextValue "E8" "  +1.78%  "

# Row 9 - USDC
Set-TextValue "E9" "  -0.48%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.678"
Set-TextValue "E10" "  +1.53%  "

# Row 11 - Dogecoin
Set-TextValue "E11" "  +8.13%  "

# Row 12 - Avalanche
Set-TextValue "D12" "56.45"
Set-TextValue "E12" "  +5.79%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  +17.02%  "

# Row 14 - Polkadot
Set-TextValue "D14" "10.11"
Set-TextValue "E14" "  +2.49%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.218.61"
Set-TextValue "E15" "  +0.70%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.641.04"
Set-TextValue "E16" "  +1.03%  "

# Row 17 - TRON
Set-TextValue "E17" "  +0.69%  "

# Row 18 - Uniswap
Set-TextValue "E18" "  +3.76%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "68.303.77"
Set-TextValue "E19" "  +2.30%  "

# Row 20 - Chainlink
Set-TextValue "D20" "18.64"
Set-TextValue "E20" "  +2.15%  "

# Row 21 - Polygon
Set-TextValue "E21" "  +3.44%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "403.25"
Set-TextValue "E22" "  +3.41%  "

# Row 23 - RenderToken
Set-TextValue "D23" "13.12"
Set-TextValue "E23" "  +28.49%  "

# Row 24 - PancakeSwap
Set-TextValue "E24" "  -0.80%  "

# Row 25 - Litecoin
Set-TextValue "E25" "  +1.79%  "

# Row 26 - ImmutableX
Set-TextValue "E26" "  +3.84%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.64"
Set-TextValue "E27" "  +4.07%  "

# Row 28 - Toncoin
Set-TextValue "E28" "  +8.02%  "

# Row 29 - LEO
Set-TextValue "E29" "  +1.36%  "

# Row 30 - NEARProtocol
Set-TextValue "D30" "8.20"
Set-TextValue "E30" "  +22.26%  "

# Row 31 - Filecoin
Set-TextValue "D31" "9.19"
Set-TextValue "E31" "  +3.30%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "31.77"
Set-TextValue "E32" "  +2.52%  "

# Row 33 - Bittensor
Set-TextValue "D33" "688.81"
Set-TextValue "E33" "  +15.76%  "

# Row 34 - Cosmos
Set-TextValue "D34" "12.26"
Set-TextValue "E34" "  +3.38%  "

# Row 35 - Hedera
Set-TextValue "E35" "  +6.00%  "

# Row 36 - OKB
Set-TextValue "D36" "64.72"
Set-TextValue "E36" "  -0.21%  "

# Row 37 - InjectiveProtocol
Set-TextValue "D37" "42.83"
Set-TextValue "E37" "  +4.46%  "

# Row 38 - TheGraph
Set-TextValue "D38" "0.425"
Set-TextValue "E38" "  +14.88%  "

# Row 39 - Dai
Set-TextValue "E39" "  -0.20%  "

# Row 40 - PEPE
Set-TextValue "D40" "0.0₃0787"
Set-TextValue "E40" "  +7.40%  "

# Row 41 - Fetch.AI
Set-TextValue "E41" "  +20.44%  "

# Row 42 - Kaspa
Set-TextValue "E42" "  +5.12%  "

# Row 43 - now ThetaToken (was Maker)
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D43" "3.15"
Set-TextValue "E43" "  +14.44%  "

# Row 44 - now Maker (was ThetaToken)
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "3.213.00"
Set-TextValue "E44" "  +16.76%  "

# Row 45 - now dogwifhat (was FirstDigitalUSD)
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D45" "3.01"
Set-TextValue "E45" "  +36.97%  "

# Row 46 - now FirstDigitalUSD (was dogwifhat)
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D46" "0.999"
Set-TextValue "E46" "  +0.13%  "

# Row 47 - VeChain
Set-TextValue "D47" "0.0420"
Set-TextValue "E47" "  +3.08%  "

# Row 48 - THORChain
Set-TextValue "D48" "8.92"
Set-TextValue "E48" "  +8.69%  "

# Row 49 - Stellar
Set-TextValue "E49" "  +2.56%  "

# Row 50 - ApeXProtocol
Set-TextValue "D50" "3.11"
Set-TextValue "E50" "  +0.99%  "

# Row 51 - WEMIXToken
Set-TextValue "E51" "  +3.68%  "
